$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the 3x3 matrix values (B1:D3)
$ws.Range("B1").Value = 0.5543
$ws.Range("C1").Value = 0.023
$ws.Range("D1").Value = -0.0888

$ws.Range("B2").Value = 0.023
$ws.Range("C2").Value = 0.5518
$ws.Range("D2").Value = 0.097345

$ws.Range("B3").Value = -0.0888
$ws.Range("C3").Value = 0.097345
$ws.Range("D3").Value = 0.042318

# Apply scientific number format to the off-diagonal cells that received it
$ws.Range("C1").NumberFormat = "0.00E+00"
$ws.Range("D1").NumberFormat = "0.00E+00"
$ws.Range("B2").NumberFormat = "0.00E+00"
$ws.Range("B3").NumberFormat = "0.00E+00"

# Update the active selection
$ws.Range("D7").Select()
